$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Widen column B on Sheet1 (target stored width 30.875; the engine quantizes
# ColumnWidth to 1/7-character pixel steps, so 30.14 is the input that lands
# on the closest achievable stored width, 30.857142857142858).
$ws1.Columns.Item(2).ColumnWidth = 30.14

# Update the remembered selection on Sheet2 (D25) without leaving it the active tab.
$ws2.Range("D25").Select()

# Make Sheet1 the active tab with D3 selected.
$ws1.Range("D3").Select()
